$wb = $excel.ActiveWorkbook

# ---- LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = "Última actualización: 05:55:25"
$ws.Range("A3").Value = "Total filas: 48"

$ws.Cells.Item(31,1).Value = "05:55:25"
$ws.Cells.Item(31,2).Value = "06:44"
$ws.Cells.Item(31,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(31,4).Value = 49
$ws.Cells.Item(31,5).Value = "LP1912"
$ws.Cells.Item(32,1).Value = "04:56:49"
$ws.Cells.Item(32,2).Value = "06:46"
$ws.Cells.Item(32,3).Value = "215C_EL PATO"
$ws.Cells.Item(32,4).Value = 110
$ws.Cells.Item(32,5).Value = "LP1912"
$ws.Cells.Item(33,1).Value = "05:26:08"
$ws.Cells.Item(33,2).Value = "06:47"
$ws.Cells.Item(33,3).Value = "215C_EL PATO"
$ws.Cells.Item(33,4).Value = 81
$ws.Cells.Item(33,5).Value = "LP1912"
$ws.Cells.Item(34,1).Value = "05:55:25"
$ws.Cells.Item(34,2).Value = "06:59"
$ws.Cells.Item(34,3).Value = "14_ABASTO"
$ws.Cells.Item(34,4).Value = 64
$ws.Cells.Item(34,5).Value = "LP1912"
$ws.Cells.Item(35,1).Value = "05:26:08"
$ws.Cells.Item(35,2).Value = "07:00"
$ws.Cells.Item(35,3).Value = "14_ABASTO"
$ws.Cells.Item(35,4).Value = 94
$ws.Cells.Item(35,5).Value = "LP1912"
$ws.Cells.Item(36,1).Value = "05:55:25"
$ws.Cells.Item(36,2).Value = "07:04"
$ws.Cells.Item(36,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(36,4).Value = 69
$ws.Cells.Item(36,5).Value = "LP1912"
$ws.Cells.Item(37,1).Value = "05:26:08"
$ws.Cells.Item(37,2).Value = "07:05"
$ws.Cells.Item(37,3).Value = "15_ABASTO"
$ws.Cells.Item(37,4).Value = 99
$ws.Cells.Item(37,5).Value = "LP1912"
$ws.Cells.Item(38,1).Value = "05:26:08"
$ws.Cells.Item(38,2).Value = "07:05"
$ws.Cells.Item(38,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(38,4).Value = 99
$ws.Cells.Item(38,5).Value = "LP1912"
$ws.Cells.Item(39,1).Value = "05:26:08"
$ws.Cells.Item(39,2).Value = "07:06"
$ws.Cells.Item(39,3).Value = "10_OLMOS"
$ws.Cells.Item(39,4).Value = 100
$ws.Cells.Item(39,5).Value = "LP1912"
$ws.Cells.Item(40,1).Value = "05:26:08"
$ws.Cells.Item(40,2).Value = "07:07"
$ws.Cells.Item(40,3).Value = "225_GOMEZ"
$ws.Cells.Item(40,4).Value = 101
$ws.Cells.Item(40,5).Value = "LP1912"
$ws.Cells.Item(41,1).Value = "05:26:08"
$ws.Cells.Item(41,2).Value = "07:11"
$ws.Cells.Item(41,3).Value = "215A_EL PATO"
$ws.Cells.Item(41,4).Value = 105
$ws.Cells.Item(41,5).Value = "LP1912"
$ws.Cells.Item(42,1).Value = "05:55:25"
$ws.Cells.Item(42,2).Value = "07:15"
$ws.Cells.Item(42,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(42,4).Value = 80
$ws.Cells.Item(42,5).Value = "LP1912"
$ws.Cells.Item(43,1).Value = "05:26:08"
$ws.Cells.Item(43,2).Value = "07:16"
$ws.Cells.Item(43,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(43,4).Value = 110
$ws.Cells.Item(43,5).Value = "LP1912"
$ws.Cells.Item(44,1).Value = "05:26:08"
$ws.Cells.Item(44,2).Value = "07:21"
$ws.Cells.Item(44,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(44,4).Value = 115
$ws.Cells.Item(44,5).Value = "LP1912"
$ws.Cells.Item(45,1).Value = "05:26:08"
$ws.Cells.Item(45,2).Value = "07:23"
$ws.Cells.Item(45,3).Value = "10_OLMOS"
$ws.Cells.Item(45,4).Value = 117
$ws.Cells.Item(45,5).Value = "LP1912"
$ws.Cells.Item(46,1).Value = "05:55:25"
$ws.Cells.Item(46,2).Value = "07:30"
$ws.Cells.Item(46,3).Value = "10_OLMOS"
$ws.Cells.Item(46,4).Value = 95
$ws.Cells.Item(46,5).Value = "LP1912"
$ws.Cells.Item(47,1).Value = "05:55:25"
$ws.Cells.Item(47,2).Value = "07:31"
$ws.Cells.Item(47,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(47,4).Value = 96
$ws.Cells.Item(47,5).Value = "LP1912"
$ws.Cells.Item(48,1).Value = "05:55:25"
$ws.Cells.Item(48,2).Value = "07:31"
$ws.Cells.Item(48,3).Value = "16_SANTA ANA"
$ws.Cells.Item(48,4).Value = 96
$ws.Cells.Item(48,5).Value = "LP1912"
$ws.Cells.Item(49,1).Value = "05:55:25"
$ws.Cells.Item(49,2).Value = "07:32"
$ws.Cells.Item(49,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(49,4).Value = 97
$ws.Cells.Item(49,5).Value = "LP1912"
$ws.Cells.Item(50,1).Value = "05:55:25"
$ws.Cells.Item(50,2).Value = "07:36"
$ws.Cells.Item(50,3).Value = "27_EL RETIRO"
$ws.Cells.Item(50,4).Value = 101
$ws.Cells.Item(50,5).Value = "LP1912"
$ws.Cells.Item(51,1).Value = "05:55:25"
$ws.Cells.Item(51,2).Value = "07:39"
$ws.Cells.Item(51,3).Value = "10_OLMOS"
$ws.Cells.Item(51,4).Value = 104
$ws.Cells.Item(51,5).Value = "LP1912"
$ws.Cells.Item(52,1).Value = "05:55:25"
$ws.Cells.Item(52,2).Value = "07:47"
$ws.Cells.Item(52,3).Value = "14_ABASTO"
$ws.Cells.Item(52,4).Value = 112
$ws.Cells.Item(52,5).Value = "LP1912"
$ws.Cells.Item(53,1).Value = "05:55:25"
$ws.Cells.Item(53,2).Value = "07:51"
$ws.Cells.Item(53,3).Value = "215D_EL PATO"
$ws.Cells.Item(53,4).Value = 116
$ws.Cells.Item(53,5).Value = "LP1912"

# ---- LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = "Última actualización: 05:55:25"
$ws.Range("A3").Value = "Total filas: 9"

$ws.Cells.Item(14,1).Value = "05:55:25"
$ws.Cells.Item(14,2).Value = "07:51"
$ws.Cells.Item(14,3).Value = "215D_EL PATO"
$ws.Cells.Item(14,4).Value = 116
$ws.Cells.Item(14,5).Value = "LP1912"

# ---- 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = "Última actualización: 05:55:25"
$ws.Range("A3").Value = "Total filas: 8"

$ws.Cells.Item(13,1).Value = "05:55:25"
$ws.Cells.Item(13,2).Value = "07:35"
$ws.Cells.Item(13,3).Value = "215A_LA PLATA"
$ws.Cells.Item(13,4).Value = 100
$ws.Cells.Item(13,5).Value = "L6173"

Write-Output "edit complete"